$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 153; this shifts existing rows 153-200 down to 154-201
$ws.Rows("153:153").Insert()

# Populate the newly inserted row 153 with the new record
$ws.Range("A153").Value = 10
$ws.Range("B153").Value = "Vega Modelo de Temuco"
$ws.Range("C153").Value = "La Araucanía"
$ws.Range("D153").Value = 45146
$ws.Range("E153").Value = 9
$ws.Range("F153").Value = 100112035
$ws.Range("G153").Value = "Bruselas (repollito)"
$ws.Range("H153").Value = "Sin especificar"
$ws.Range("I153").Value = "Primera"
$ws.Range("J153").Value = 55
$ws.Range("K153").Value = 25000
$ws.Range("L153").Value = 25000
$ws.Range("M153").Value = 25000
$ws.Range("N153").Value = "$/malla 15 kilos"
$ws.Range("O153").Value = "Provincia de Quillota"
$ws.Range("P153").Value = 1667
$ws.Range("Q153").Value = 15
$ws.Range("R153").Value = "Hortaliza"
